# Applies the "Act graficos y tablas web pob" revision to indicadores/tablas/pob/330112.xlsx
#   - renames the two worksheets (Datos -> Data, Ficha técnica -> Metadata)
#   - refreshes the Data sheet's year/value series (2015-2021) on sheet "Data"
#   - rewrites the Metadata sheet's indicator documentation (lower-case keys,
#     refreshed "AFAM-PE" wording, new "observaciones" + citation rows)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Worksheet names
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsMeta = $wb.Worksheets.Item(2)
$wsData.Name = "Data"
$wsMeta.Name = "Metadata"

# ---------------------------------------------------------------------------
# 2. "Data" sheet - Fecha / Valor series
# ---------------------------------------------------------------------------
# Header row is unchanged: A1=Fecha, B1=Valor

function Set-YearCell($ws, $row, $year, $value) {
    # Years must stay text (matching the rest of the column) rather than be
    # auto-coerced to numbers, so force a text format before writing.
    $cell = $ws.Range("A" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $year
    $ws.Range("B" + $row).Value = $value
}

Set-YearCell $wsData 2 "2021" 17.4
Set-YearCell $wsData 3 "2020" 17.6
Set-YearCell $wsData 4 "2019" 13.9
Set-YearCell $wsData 5 "2018" 16
Set-YearCell $wsData 6 "2017" 14
Set-YearCell $wsData 7 "2016" 15.6
Set-YearCell $wsData 8 "2015" 17.3

# ---------------------------------------------------------------------------
# 3. "Metadata" sheet - indicator documentation
# ---------------------------------------------------------------------------
# Row 1 (blank header) is unchanged.

$wsMeta.Range("A2").Value = "nomindicador"
$wsMeta.Range("B2").Value = "Porcentaje de NNA en hogares pobres que no perciben AFAM-PE"

$wsMeta.Range("A3").Value = "derecho"
$wsMeta.Range("B3").Value = "Seguridad Social"

$wsMeta.Range("A4").Value = "conindicador"
$wsMeta.Range("B4").Value = "NNA en hogares pobres que no perciben AFAM-PE"

$wsMeta.Range("A5").Value = "tipoind"
$wsMeta.Range("B5").Value = "Resultados"

$wsMeta.Range("A6").Value = "definicion"
$wsMeta.Range("B6").Value = "El indicador mide el porcentaje de niños, niñas y adolescentes (hasta 17 años de edad) que residen en hogares en los que ningún integrante percibe Asignaciones Familiares - Plan de Equidad en el total de hogares pobres."

$wsMeta.Range("A7").Value = "calculo"
$wsMeta.Range("B7").Value = "Para cada año calcular: (Cantidad de niños, niñas y adolescentes (hasta 17 años de edad) que residen en hogares en los que ninguno de sus integrantes es beneficiario de Asignaciones Familiares - Plan de Equidad/Cantidad de niños, niñas y adolescentes que residen en hogares pobres)*100"

$wsMeta.Range("A8").Value = "observaciones"
$wsMeta.Range("B8").Value = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. `nEn julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH.  `nLos indicadores de trabajo y seguridad social del año 2020 se construyen con la encuesta presencial realizada hasta marzo de 2020 y posteriormente con la encuesta telefónica panel (siempre que la información haya sido incluida en el formulario). Para el 2021, se calculan a partir de la encuesta telefónica del primer semestre de 2021 y el formulario telefónico de modalidad panel del segundo semestre de 2021. En el segundo semestre de 2021 el quintil de ingresos del hogar corresponde a los ingresos declarados durante la implantación del panel en la encuesta presencial."

$wsMeta.Range("A9").Value = "cita"
$wsMeta.Range("B9").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE"

$wsMeta.Range("A10").Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$wsMeta.Range("B10").Value = " "
